$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Compartments sheet: insert a new "Type" column between the existing
#    "Initial volume" (D) and "Comments" (old E, now F) columns.
# ---------------------------------------------------------------------
$compartments = $wb.Worksheets.Item("Compartments")

# Shift columns E:F -> F:G and create a new, empty column E.
$compartments.Range("E1").EntireColumn.Insert()

# Header for the newly inserted column.
$compartments.Range("E1").Value = "Type"

# Re-establish the AutoFilter over the widened data range (A1:F3).
$compartments.AutoFilterMode = $false
[void]$compartments.Range("A1:F3").AutoFilter()

# ---------------------------------------------------------------------
# 2. Fix up the two "_FilterDatabase" defined names that point at the
#    Compartments sheet so they cover the extra column too.
# ---------------------------------------------------------------------
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "Compartments!_FilterDatabase" -or $n.Name -eq "Compartments!_FilterDatabase_0") {
        $n.RefersTo = "=Compartments!`$A`$1:`$F`$3"
    }
}

# ---------------------------------------------------------------------
# 3. Reactions sheet: rows 2, 4, 5 and 6 go back to the default row
#    height (no more explicit 26pt height).
# ---------------------------------------------------------------------
$reactions = $wb.Worksheets.Item("Reactions")
$reactions.Rows.Item(2).AutoFit()
$reactions.Rows.Item(4).AutoFit()
$reactions.Rows.Item(5).AutoFit()
$reactions.Rows.Item(6).AutoFit()

# ---------------------------------------------------------------------
# 4. Make "Compartments" the active/selected sheet (was "Reactions").
# ---------------------------------------------------------------------
$compartments.Activate()
